$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Spherical")
$ws.Range("A2").Value = [double]"0.6675858082427026"
$ws.Range("B2").Value = [double]"0.1973354237725259"
$ws.Range("C2").Value = [double]"0.2781024412082834"
$ws.Range("D2").Value = [double]"0.4442245195534865"
$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"6.979693377975784e-29"
$ws.Range("G2").Value = [double]"4.406330720924272e-15"
$ws.Range("H2").Value = [double]"8.354455923622905e-15"
$ws.Range("I2").Value = [double]"0.6513796849758856"
$ws.Range("J2").Value = [double]"0.2442157708210946"
$ws.Range("K2").Value = [double]"0.4838358795907293"
$ws.Range("L2").Value = [double]"0.4941819207752289"
$ws.Range("M2").Value = [double]"0.7443853761732526"
$ws.Range("N2").Value = [double]"0.1251941513667706"
$ws.Range("O2").Value = [double]"0.1597491092106763"
$ws.Range("P2").Value = [double]"0.3538278555551705"

$ws = $wb.Worksheets.Item("Gaussian")
$ws.Range("A2").Value = [double]"0.7082707905677209"
$ws.Range("B2").Value = [double]"0.1731830607646702"
$ws.Range("C2").Value = [double]"0.2437041923490998"
$ws.Range("D2").Value = [double]"0.4161526892435878"
$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"7.036755824661521e-24"
$ws.Range("G2").Value = [double]"1.431912009762658e-12"
$ws.Range("H2").Value = [double]"2.652688414545048e-12"
$ws.Range("I2").Value = [double]"0.6686518898700511"
$ws.Range("J2").Value = [double]"0.2321162325835804"
$ws.Range("K2").Value = [double]"0.4502923171128647"
$ws.Range("L2").Value = [double]"0.481784425426539"
$ws.Range("M2").Value = [double]"0.7946388445521382"
$ws.Range("N2").Value = [double]"0.1005811607923515"
$ws.Range("O2").Value = [double]"0.1272597757958034"
$ws.Range("P2").Value = [double]"0.3171453307118859"

$ws = $wb.Worksheets.Item("Exponential")
$ws.Range("A2").Value = [double]"0.6292361889479723"
$ws.Range("B2").Value = [double]"0.2201014144031727"
$ws.Range("C2").Value = [double]"0.3093648956339847"
$ws.Range("D2").Value = [double]"0.4691496716434667"
$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"1.233457981022916e-29"
$ws.Range("G2").Value = [double]"1.525742510176214e-15"
$ws.Range("H2").Value = [double]"3.512062045327384e-15"
$ws.Range("I2").Value = [double]"0.5831568642345634"
$ws.Range("J2").Value = [double]"0.292007273601932"
$ws.Range("K2").Value = [double]"0.5586679374723076"
$ws.Range("L2").Value = [double]"0.5403769736044755"
$ws.Range("M2").Value = [double]"0.7132259063696762"
$ws.Range("N2").Value = [double]"0.1404553415158181"
$ws.Range("O2").Value = [double]"0.1803210790901188"
$ws.Range("P2").Value = [double]"0.3747737204178251"

$ws = $wb.Worksheets.Item("Linear")
$ws.Range("A2").Value = [double]"-0.05017696330948196"
$ws.Range("B2").Value = [double]"0.6234304107031905"
$ws.Range("C2").Value = [double]"0.5921680198200978"
$ws.Range("D2").Value = [double]"0.7895760955748283"
$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"2.194841122755544e-30"
$ws.Range("G2").Value = [double]"6.987666884789323e-16"
$ws.Range("H2").Value = [double]"1.481499619559703e-15"
$ws.Range("I2").Value = [double]"-0.08863502004839785"
$ws.Range("J2").Value = [double]"0.7626114403159991"
$ws.Range("K2").Value = [double]"0.9501502032838559"
$ws.Range("L2").Value = [double]"0.8732762680366386"
$ws.Range("M2").Value = [double]"-0.01636870287358216"
$ws.Range("N2").Value = [double]"0.4977939654902045"
$ws.Range("O2").Value = [double]"0.3979895811419075"
$ws.Range("P2").Value = [double]"0.7055451548201607"

$ws = $wb.Worksheets.Item("Power")
$ws.Range("A2").Value = [double]"0.3200064928743519"
$ws.Range("B2").Value = [double]"0.4036735200198026"
$ws.Range("C2").Value = [double]"0.4646543711342551"
$ws.Range("D2").Value = [double]"0.6353530672152317"
$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"2.576123893612367e-30"
$ws.Range("G2").Value = [double]"7.892850043672255e-16"
$ws.Range("H2").Value = [double]"1.605030807683257e-15"
$ws.Range("I2").Value = [double]"0.2840702064200564"
$ws.Range("J2").Value = [double]"0.5015236888327036"
$ws.Range("K2").Value = [double]"0.770008918542557"
$ws.Range("L2").Value = [double]"0.7081833723215363"
$ws.Range("M2").Value = [double]"0.3801382620230105"
$ws.Range("N2").Value = [double]"0.3035939927418204"
$ws.Range("O2").Value = [double]"0.3031425412414653"
$ws.Range("P2").Value = [double]"0.5509936412898251"

$ws = $wb.Worksheets.Item("HoleEffect")
$ws.Range("A2").Value = [double]"0.6416256090604304"
$ws.Range("B2").Value = [double]"0.2127465194293359"
$ws.Range("C2").Value = [double]"0.3042790984108429"
$ws.Range("D2").Value = [double]"0.4612445332243363"
$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"2.74346923043347e-29"
$ws.Range("G2").Value = [double]"2.418830424768703e-15"
$ws.Range("H2").Value = [double]"5.237813695076859e-15"
$ws.Range("I2").Value = [double]"0.5855851084770898"
$ws.Range("J2").Value = [double]"0.2903062380802658"
$ws.Range("K2").Value = [double]"0.5498088035944768"
$ws.Range("L2").Value = [double]"0.5388007406084978"
$ws.Range("M2").Value = [double]"0.7270330263336232"
$ws.Range("N2").Value = [double]"0.1336929323827744"
$ws.Range("O2").Value = [double]"0.1757018704271285"
$ws.Range("P2").Value = [double]"0.3656404413939662"

